# Applies the rests.docx implementation-notes clarifications:
#  - removes stale proofErr (spelling/grammar) markers
#  - adds a trailing period + splits the second paragraph
#  - appends a new paragraph about augmentationDot
#  - re-expresses the header's STYLEREF field as a complex field (fldChar triplet)

$d = $word.ActiveDocument

# --- Body content -----------------------------------------------------
$bodyXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="705B36CF" w14:textId="77777777" w:rsidR="0002482F" w:rsidRDefault="0002482F" w:rsidP="0002482F"><w:pPr><w:pStyle w:val="Body"/></w:pPr><w:r><w:t xml:space="preserve">Scoring applications should draw multiple measure rests using primitives to provide variable width and line thickness rather than using </w:t></w:r><w:r w:rsidRPr="007F76C0"><w:rPr><w:rFonts w:ascii="Avenir Heavy" w:hAnsi="Avenir Heavy"/></w:rPr><w:t>restHBar</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p w14:paraId="12924381" w14:textId="5748A6B9" w:rsidR="0095527F" w:rsidRPr="003D7783" w:rsidRDefault="0002482F" w:rsidP="0002482F"><w:pPr><w:pStyle w:val="Body"/></w:pPr><w:r><w:t xml:space="preserve">“Old style” multiple measure rests can be created by laying out </w:t></w:r><w:r w:rsidRPr="007F76C0"><w:rPr><w:rFonts w:ascii="Avenir Heavy" w:hAnsi="Avenir Heavy"/></w:rPr><w:t>restLonga</w:t></w:r><w:r><w:t xml:space="preserve"> (four bars), </w:t></w:r><w:r w:rsidRPr="007F76C0"><w:rPr><w:rFonts w:ascii="Avenir Heavy" w:hAnsi="Avenir Heavy"/></w:rPr><w:t>restDoubleWhole</w:t></w:r><w:r><w:t xml:space="preserve"> (two bars) and </w:t></w:r><w:r w:rsidRPr="007F76C0"><w:rPr><w:rFonts w:ascii="Avenir Heavy" w:hAnsi="Avenir Heavy"/></w:rPr><w:t>restWhole</w:t></w:r><w:r><w:t xml:space="preserve"> (one bar) next to each other.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/></w:pPr><w:r><w:t xml:space="preserve">For dotted rests, the augmentation dot glyph </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Avenir Heavy" w:hAnsi="Avenir Heavy"/></w:rPr><w:t>augmentationDot</w:t></w:r><w:r><w:t xml:space="preserve"> should be used.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Content.InsertXML($bodyXml)

# --- Header field -------------------------------------------------------
$headerXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="3284E791" w14:textId="77777777" w:rsidR="0052215B" w:rsidRPr="00FC768B" w:rsidRDefault="00194177" w:rsidP="00FC768B"><w:pPr><w:pStyle w:val="Header"/></w:pPr><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> STYLEREF "Heading 1" \* MERGEFORMAT </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:t>Staff brackets (U+E000–U+E01F)</w:t></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
$hdr.Range.InsertXML($headerXml)

Write-Output "done"
